# Fix the typo "Ore and mineral indutries" -> "Ore and mineral industries"
# wherever it appears in the active worksheet (column B, CESI Sector values).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$old = "Ore and mineral indutries"
$new = "Ore and mineral industries"

$used = $ws.UsedRange
$rows = $used.Rows.Count
$cols = $used.Columns.Count

for ($r = 1; $r -le $rows; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        if ($cell.Value2 -eq $old) {
            $cell.Value = $new
        }
    }
}
